$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -162.4

$ws.Range("B3").Value = -428.4
$ws.Range("C3").Value = -324.5

$ws.Range("C4").Value = -296.9

$ws.Range("C5").Value = -96.7

$ws.Range("C6").Value = -63.1

$ws.Range("C7").Value = -147.9

$ws.Range("C8").Value = 6.7

$ws.Range("C10").Value = -7.3

$ws.Range("C11").Value = -60

$ws.Range("C12").Value = -155.9

$ws.Range("C13").Value = 131.4

$ws.Range("C14").Value = 206.3

$ws.Range("C16").Value = -90.3

$ws.Range("C17").Value = 183.8

$ws.Range("C18").Value = -13.6

$ws.Range("C19").Value = -261.8

$ws.Range("C21").Value = 59.5

$ws.Range("C23").Value = 27.7

$ws.Range("C24").Value = 131.9
